# Fix UI issues: remove duplicate Offer Details and enable Add New Candidate button
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dash = [string][char]0x2014

# Row 16
$ws.Range("M16").Value = "Offer"
$ws.Range("N16").Value = "In Notice"

# Row 17
$ws.Range("O17").Value = $dash
$ws.Range("P17").Value = $dash
$ws.Range("Q17").Value = ""
$ws.Range("V17").Value = $dash
$ws.Range("Y17").Value = $dash
$ws.Range("AF17").Value = $dash

# Row 18
$ws.Range("AE18").Value = "No"

# Row 19
$ws.Range("A19").Value = "Invalid Date"
$ws.Range("L19").Value = ""
$ws.Range("O19").Value = $dash
$ws.Range("P19").Value = $dash
$ws.Range("Q19").Value = ""
$ws.Range("W19").Value = $dash
$ws.Range("Y19").Value = $dash
$ws.Range("AF19").Value = $dash
